$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.070.57"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.375.71"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.60"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.46"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.71"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0790"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.49"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.741.03"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.376.90"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.114.64"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.96"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.11"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.98"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.52"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.97"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("E32").Value = "  +12.67%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.07"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.93"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0738"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "129.35"
$ws.Range("E36").Value = "  +12.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.85"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.32"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.06"
$ws.Range("E42").Value = "  -6.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.930.42"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.27"
$ws.Range("E47").Value = "  -7.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.600.33"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.98"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.64"
$ws.Range("E51").Value = "  -0.98%  "
